$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -11
$ws.Range("F3").Value = -3
$ws.Range("F4").Value = -3
$ws.Range("F5").Value = 1
$ws.Range("F18").Value = -1
$ws.Range("F21").Value = -1
